$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Inscritos" (column E) values
$ws.Range("E12").Value = 34
$ws.Range("E15").Value = 112
$ws.Range("E16").Value = 333
$ws.Range("E17").Value = 30

# Update "Pagos" (column F) and "Inscrições homologadas" (column H) for row 18
$ws.Range("F18").Value = 34
$ws.Range("H18").Value = 57
